$wb = $excel.ActiveWorkbook

# Remember which sheet is currently active so adding the new sheet doesn't
# change the workbook's active-tab selection.
$prevActive = $wb.ActiveSheet

# Add the new day sheet ("2024-08-22") right after the last existing sheet
# ("2024-08-21"), matching the existing sheet naming/ordering convention.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2024-08-22"

# Populate it the same way every other archived-tasks day sheet is populated.
$ws.Range("A1").Value = "Archived Tasks"
$ws.Range("A2").Value = "1H Py Apps Course"

# Copy the header cell's formatting (bold + box border) from the previous
# day's sheet so the new header looks identical to the others.
$prevSheet = $wb.Worksheets.Item("2024-08-21")
$prevSheet.Range("A1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Restore the original active sheet/selection.
$prevActive.Activate()
